$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Column B was "Product-Code" (e.g. "VG-WHITE"); it now also encodes a
# size segment, so the header and existing codes are updated to the
# new "Product-Size-Color" scheme (e.g. "VG-XS-WHITE").
$ws.Range("B1").Value = "Product-Size-Color"
$ws.Range("B2").Value = "VG-XS-WHITE"
$ws.Range("B4").Value = "BN-XS-PURPLE"
$ws.Range("B5").Value = "BN-S-RED"

# Row 7 previously duplicated the GH-BLUE row; it now represents a
# distinct red-colored variant (new shared string "GH-RED").
$ws.Range("B7").Value = "GH-RED"

# Move the active selection to B8, matching the updated sheet view.
$ws.Range("B8").Select()
